# Applies the "Update MDSC and BDESC" commit to the workbook.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) About sheet: add explanatory notes about the 50/50 commercial /
#    residential rooftop-solar split, with a source reference.
# ---------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")
$about.Range("C23").Value = "We assume a 50/50 split between commercial and residential rooftop solar."
$about.Range("C24").Value = "See Solar Power Europe's Figure 11:"
$about.Range("C25").Value = "https://www.solarpowereurope.org/insights/outlooks/eu-market-outlook-for-solar-power-2023-2027/detail"

# ---------------------------------------------------------------------
# 2) BDESC-urban-residential: the rooftop-solar capacity pulled from the
#    Raw data sheet is now halved (the other half now goes to
#    BDESC-commercial, see below).
# ---------------------------------------------------------------------
$urban = $wb.Worksheets.Item("BDESC-urban-residential")
$urban.Range("B8").Formula = "='Raw data'!B21/2"
$urban.Range("C8").Formula = "='Raw data'!C21/2"
$urban.Range("D8").Formula = "='Raw data'!D21/2"
$urban.Range("G8").Formula = "='Raw data'!E21/2"
$urban.Range("L8").Formula = "='Raw data'!F21/2"
$urban.Range("Q8").Formula = "='Raw data'!G21/2"
$urban.Range("V8").Formula = "='Raw data'!H21/2"
$urban.Range("AA8").Formula = "='Raw data'!I21/2"
$urban.Range("AF8").Formula = "='Raw data'!J21/2"

# ---------------------------------------------------------------------
# 3) BDESC-commercial: now takes the other half of the rooftop-solar
#    capacity, pulled directly from BDESC-urban-residential row 8.
# ---------------------------------------------------------------------
$commercial = $wb.Worksheets.Item("BDESC-commercial")
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF")
foreach ($col in $cols) {
    $addr = $col + "8"
    $commercial.Range($addr).Formula = "='BDESC-urban-residential'!" + $addr
}

# ---------------------------------------------------------------------
# 4) Selections / active-sheet bookkeeping, matching the author's last
#    on-screen state when they saved the workbook.
# ---------------------------------------------------------------------
$raw = $wb.Worksheets.Item("Raw data")
$raw.Activate() | Out-Null
$raw.Range("F19").Select() | Out-Null

$commercial.Activate() | Out-Null
$commercial.Range("B8:AF8").Select() | Out-Null

$about.Activate() | Out-Null
$about.Range("C26").Select() | Out-Null

$urban.Activate() | Out-Null
$urban.Range("AF8").Select() | Out-Null

$wb.Application.Calculate() | Out-Null
